# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to the freshly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 63
$ws1.Range("F9").Value  = 868
$ws1.Range("F12").Value = 10307
$ws1.Range("F14").Value = 273
$ws1.Range("F15").Value = 17
$ws1.Range("F16").Value = 649
$ws1.Range("G16").Value = 59.9
$ws1.Range("F17").Value = 11875
$ws1.Range("F18").Value = 12260
$ws1.Range("F20").Value = 105

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 63
$ws4.Range("F10").Value = 869
$ws4.Range("F13").Value = 10307
$ws4.Range("F15").Value = 273
$ws4.Range("F16").Value = 17
$ws4.Range("F17").Value = 649
$ws4.Range("G17").Value = 59.9
$ws4.Range("F18").Value = 11875
$ws4.Range("F19").Value = 12260
$ws4.Range("F21").Value = 105
